$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 142.2381
$ws.Range("I33").Value = 158.52942
$ws.Range("J33").Value = 73
$ws.Range("K33").Value = 158.52942
$ws.Range("L33").Value = 73
$ws.Range("M33").Value = 70.47058000000001
$ws.Range("N33").Value = -531
# Row 40
$ws.Range("H40").Value = 1428.6111
$ws.Range("I40").Value = 1352.5
$ws.Range("J40").Value = 1523.75
$ws.Range("K40").Value = 1352.5
$ws.Range("L40").Value = 1523.75
$ws.Range("M40").Value = -1177.5
$ws.Range("N40").Value = -1873.75
# Row 64
$ws.Range("H64").Value = 3691.9666
$ws.Range("I64").Value = 3651
$ws.Range("J64").Value = 3709.524
$ws.Range("K64").Value = 3651
$ws.Range("L64").Value = 3709.524
$ws.Range("M64").Value = -3403
$ws.Range("N64").Value = -4205.523999999999
# Row 67
$ws.Range("H67").Value = 3691.9666
$ws.Range("I67").Value = 3651
$ws.Range("J67").Value = 3709.524
$ws.Range("K67").Value = 3651
$ws.Range("L67").Value = 3709.524
$ws.Range("M67").Value = -2793
$ws.Range("N67").Value = -5425.523999999999
# Row 98
$ws.Range("H98").Value = 9383.625
$ws.Range("I98").Value = 6113.8
$ws.Range("J98").Value = 14833.333
$ws.Range("K98").Value = 6113.8
$ws.Range("L98").Value = 14833.333
$ws.Range("M98").Value = -4615.8
$ws.Range("N98").Value = -17829.333
# Row 122
$ws.Range("H122").Value = 9383.625
$ws.Range("I122").Value = 6113.8
$ws.Range("J122").Value = 14833.333
$ws.Range("K122").Value = 18341.4
$ws.Range("L122").Value = 44499.999
$ws.Range("M122").Value = -15891.4
$ws.Range("N122").Value = -49399.999
# Row 137
$ws.Range("H137").Value = 2662.8696
$ws.Range("I137").Value = 1647.8
$ws.Range("K137").Value = 4943.4
$ws.Range("M137").Value = -2393.4

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 1828.5
$ws.Range("I5").Value = 2154.6
$ws.Range("J5").Value = 198
$ws.Range("K5").Value = 2154.6
$ws.Range("L5").Value = 198
$ws.Range("M5").Value = -2042.6
$ws.Range("N5").Value = -422
# Row 32
$ws.Range("H32").Value = 12157.238
$ws.Range("I32").Value = 13743.529
$ws.Range("J32").Value = 5415.5
$ws.Range("K32").Value = 13743.529
$ws.Range("L32").Value = 5415.5
$ws.Range("M32").Value = -13456.529
$ws.Range("N32").Value = -5989.5
# Row 74
$ws.Range("H74").Value = 1357.8422
$ws.Range("I74").Value = 1185.9773
$ws.Range("J74").Value = 1939.5385
$ws.Range("K74").Value = 1185.9773
$ws.Range("L74").Value = 1939.5385
$ws.Range("M74").Value = -311.9773
$ws.Range("N74").Value = -3687.5385
# Row 77
$ws.Range("H77").Value = 1357.8422
$ws.Range("I77").Value = 1185.9773
$ws.Range("J77").Value = 1939.5385
$ws.Range("K77").Value = 5929.886500000001
$ws.Range("L77").Value = 9697.692500000001
$ws.Range("M77").Value = -1561.886500000001
$ws.Range("N77").Value = -18433.6925
# Row 102
$ws.Range("H102").Value = 93257.18
$ws.Range("I102").Value = 2528.75
$ws.Range("J102").Value = 335199.66
$ws.Range("K102").Value = 2528.75
$ws.Range("L102").Value = 335199.66
$ws.Range("M102").Value = -906.75
$ws.Range("N102").Value = -338443.66
# Row 119
$ws.Range("H119").Value = 43666.5
$ws.Range("J119").Value = 43666.5
$ws.Range("L119").Value = 43666.5
$ws.Range("N119").Value = -53342.5
# Row 122
$ws.Range("H122").Value = 4367.4165
$ws.Range("I122").Value = 5066.4443
$ws.Range("J122").Value = 2270.3333
$ws.Range("K122").Value = 15199.3329
$ws.Range("L122").Value = 6810.999899999999
$ws.Range("M122").Value = -12749.3329
$ws.Range("N122").Value = -11710.9999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 1828.5
$ws.Range("I4").Value = 2154.6
$ws.Range("J4").Value = 198
$ws.Range("K4").Value = 2154.6
$ws.Range("L4").Value = 198
$ws.Range("M4").Value = -2039.6
$ws.Range("N4").Value = -428

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 52542.1
$ws.Range("I62").Value = 85684.164
$ws.Range("J62").Value = 2829
$ws.Range("K62").Value = 85684.164
$ws.Range("L62").Value = 2829
$ws.Range("M62").Value = -85060.164
$ws.Range("N62").Value = -4077
# Row 65
$ws.Range("H65").Value = 52542.1
$ws.Range("I65").Value = 85684.164
$ws.Range("J65").Value = 2829
$ws.Range("K65").Value = 428420.82
$ws.Range("L65").Value = 14145
$ws.Range("M65").Value = -425300.82
$ws.Range("N65").Value = -20385

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 3000
$ws.Range("J17").Value = 4000
$ws.Range("L17").Value = 12000
$ws.Range("N17").Value = -12338
# Row 34
$ws.Range("H34").Value = 1064.125
$ws.Range("I34").Value = 418.33334
$ws.Range("J34").Value = 3001.5
$ws.Range("K34").Value = 1255.00002
$ws.Range("L34").Value = 9004.5
$ws.Range("M34").Value = -1171.00002
$ws.Range("N34").Value = -9172.5
# Row 39
$ws.Range("H39").Value = 5173.25
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 5598.091
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 16794.273
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -17382.273
# Row 55
$ws.Range("H55").Value = 4672.2666
$ws.Range("I55").Value = 2899
$ws.Range("J55").Value = 4798.9287
$ws.Range("K55").Value = 8697
$ws.Range("L55").Value = 14396.7861
$ws.Range("M55").Value = -8520
$ws.Range("N55").Value = -14750.7861
# Row 68
$ws.Range("H68").Value = 159455.66
$ws.Range("I68").Value = 238648.53
$ws.Range("J68").Value = 1069.9048
$ws.Range("K68").Value = 715945.59
$ws.Range("L68").Value = 3209.7144
$ws.Range("M68").Value = -715134.59
$ws.Range("N68").Value = -4831.7144
# Row 71
$ws.Range("H71").Value = 159455.66
$ws.Range("I71").Value = 238648.53
$ws.Range("J71").Value = 1069.9048
$ws.Range("K71").Value = 2147836.77
$ws.Range("L71").Value = 9629.1432
$ws.Range("M71").Value = -2143780.77
$ws.Range("N71").Value = -17741.1432
# Row 131
$ws.Range("H131").Value = 2657.2808
$ws.Range("I131").Value = 478.94116
$ws.Range("J131").Value = 3583.075
$ws.Range("K131").Value = 1436.82348
$ws.Range("L131").Value = 10749.225
$ws.Range("M131").Value = 3603.17652
$ws.Range("N131").Value = -20829.225

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 28004
$ws.Range("J5").Value = 30005
$ws.Range("L5").Value = 30005
$ws.Range("N5").Value = -30229

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 912.6667
$ws.Range("I46").Value = 675
$ws.Range("J46").Value = 1388
$ws.Range("K46").Value = 675
$ws.Range("L46").Value = 1388
$ws.Range("M46").Value = -487
$ws.Range("N46").Value = -1764

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 68
$ws.Range("H68").Value = 46566.668
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 46566.668
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 46566.668
$ws.Range("N68").Value = -48188.668
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 46566.668
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 46566.668
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 139700.004
$ws.Range("N71").Value = -147812.004
$ws.Range("M71").ClearContents()
# Row 101
$ws.Range("H101").Value = 43999.5
$ws.Range("J101").Value = 43999.5
$ws.Range("L101").Value = 43999.5
$ws.Range("N101").Value = -50489.5
